$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Existing-cell edits ---
$ws.Range("N774").Value = 2

# Q779:Q783 change from blank/inlineStr to numeric 0
$ws.Range("Q779").Value = 0
$ws.Range("Q780").Value = 0
$ws.Range("Q781").Value = 0
$ws.Range("Q782").Value = 0
$ws.Range("Q783").Value = 0

# --- New rows 784:799 ---
# Match the date-time number format used by the rest of column A
$dateFmt = $ws.Range("A783").NumberFormat

# Row 784
$ws.Range("A784").NumberFormat = $dateFmt
$ws.Range("A784").Value = 45666
$ws.Range("B784").Value = 169.1000061035156
$ws.Range("C784").Value = 173.4299926757812
$ws.Range("D784").Value = 168.3000030517578
$ws.Range("E784").Value = 170.1999969482422
$ws.Range("F784").Value = 4982727
$ws.Range("G784").Value = 2025
$ws.Range("H784").Value = 1
$ws.Range("I784").Value = 9
$ws.Range("J784").Value = 0
$ws.Range("K784").Value = 0
$ws.Range("L784").Value = 0
$ws.Range("M784").Value = 2
$ws.Range("N784").Value = 0
$ws.Range("O784").Value = 0
$ws.Range("P784").Value = 0

# Row 785
$ws.Range("A785").NumberFormat = $dateFmt
$ws.Range("A785").Value = 45667
$ws.Range("B785").Value = 171.4100036621094
$ws.Range("C785").Value = 171.4100036621094
$ws.Range("D785").Value = 168.4499969482422
$ws.Range("E785").Value = 169.6999969482422
$ws.Range("F785").Value = 3475881
$ws.Range("G785").Value = 2025
$ws.Range("H785").Value = 1
$ws.Range("I785").Value = 10
$ws.Range("J785").Value = 0
$ws.Range("K785").Value = 0
$ws.Range("L785").Value = 0
$ws.Range("M785").Value = 2
$ws.Range("N785").Value = 0
$ws.Range("O785").Value = 1
$ws.Range("P785").Value = 1

# Row 786
$ws.Range("A786").NumberFormat = $dateFmt
$ws.Range("A786").Value = 45670
$ws.Range("B786").Value = 167.9900054931641
$ws.Range("C786").Value = 168
$ws.Range("D786").Value = 161.6600036621094
$ws.Range("E786").Value = 164.5399932861328
$ws.Range("F786").Value = 9025308
$ws.Range("G786").Value = 2025
$ws.Range("H786").Value = 1
$ws.Range("I786").Value = 13
$ws.Range("J786").Value = 0
$ws.Range("K786").Value = 0
$ws.Range("L786").Value = 0
$ws.Range("M786").Value = 3
$ws.Range("N786").Value = 0
$ws.Range("O786").Value = 0
$ws.Range("P786").Value = 0

# Row 787
$ws.Range("A787").NumberFormat = $dateFmt
$ws.Range("A787").Value = 45671
$ws.Range("B787").Value = 164.0500030517578
$ws.Range("C787").Value = 166.1900024414062
$ws.Range("D787").Value = 163.1199951171875
$ws.Range("E787").Value = 164.9299926757812
$ws.Range("F787").Value = 15903730
$ws.Range("G787").Value = 2025
$ws.Range("H787").Value = 1
$ws.Range("I787").Value = 14
$ws.Range("J787").Value = 0
$ws.Range("K787").Value = 0
$ws.Range("L787").Value = 0
$ws.Range("M787").Value = 3
$ws.Range("N787").Value = 0
$ws.Range("O787").Value = 0
$ws.Range("P787").Value = 0

# Row 788
$ws.Range("A788").NumberFormat = $dateFmt
$ws.Range("A788").Value = 45672
$ws.Range("B788").Value = 167.5800018310547
$ws.Range("C788").Value = 173.3800048828125
$ws.Range("D788").Value = 165.6399993896484
$ws.Range("E788").Value = 172.7700042724609
$ws.Range("F788").Value = 7801286
$ws.Range("G788").Value = 2025
$ws.Range("H788").Value = 1
$ws.Range("I788").Value = 15
$ws.Range("J788").Value = 0
$ws.Range("K788").Value = 0
$ws.Range("L788").Value = 0
$ws.Range("M788").Value = 3
$ws.Range("N788").Value = 0
$ws.Range("O788").Value = 0
$ws.Range("P788").Value = 0

# Row 789
$ws.Range("A789").NumberFormat = $dateFmt
$ws.Range("A789").Value = 45673
$ws.Range("B789").Value = 172.5899963378906
$ws.Range("C789").Value = 176.8500061035156
$ws.Range("D789").Value = 171.6100006103516
$ws.Range("E789").Value = 174.9799957275391
$ws.Range("F789").Value = 9010520
$ws.Range("G789").Value = 2025
$ws.Range("H789").Value = 1
$ws.Range("I789").Value = 16
$ws.Range("J789").Value = 0
$ws.Range("K789").Value = 0
$ws.Range("L789").Value = 0
$ws.Range("M789").Value = 3
$ws.Range("N789").Value = 1
$ws.Range("O789").Value = 0
$ws.Range("P789").Value = 0

# Row 790
$ws.Range("A790").NumberFormat = $dateFmt
$ws.Range("A790").Value = 45674
$ws.Range("B790").Value = 173.4199981689453
$ws.Range("C790").Value = 175.5
$ws.Range("D790").Value = 172.1199951171875
$ws.Range("E790").Value = 172.8800048828125
$ws.Range("F790").Value = 3549636
$ws.Range("G790").Value = 2025
$ws.Range("H790").Value = 1
$ws.Range("I790").Value = 17
$ws.Range("J790").Value = 0
$ws.Range("K790").Value = 0
$ws.Range("L790").Value = 0
$ws.Range("M790").Value = 3
$ws.Range("N790").Value = 0
$ws.Range("O790").Value = 0
$ws.Range("P790").Value = 0

# Row 791
$ws.Range("A791").NumberFormat = $dateFmt
$ws.Range("A791").Value = 45677
$ws.Range("B791").Value = 173.4400024414062
$ws.Range("C791").Value = 174.2200012207031
$ws.Range("D791").Value = 166.5
$ws.Range("E791").Value = 172.1900024414062
$ws.Range("F791").Value = 5955385
$ws.Range("G791").Value = 2025
$ws.Range("H791").Value = 1
$ws.Range("I791").Value = 20
$ws.Range("J791").Value = 0
$ws.Range("K791").Value = 0
$ws.Range("L791").Value = 0
$ws.Range("M791").Value = 4
$ws.Range("N791").Value = 0
$ws.Range("O791").Value = 0
$ws.Range("P791").Value = 0

# Row 792
$ws.Range("A792").NumberFormat = $dateFmt
$ws.Range("A792").Value = 45678
$ws.Range("B792").Value = 171.75
$ws.Range("C792").Value = 172.0899963378906
$ws.Range("D792").Value = 166.1399993896484
$ws.Range("E792").Value = 168.1999969482422
$ws.Range("F792").Value = 4526759
$ws.Range("G792").Value = 2025
$ws.Range("H792").Value = 1
$ws.Range("I792").Value = 21
$ws.Range("J792").Value = 0
$ws.Range("K792").Value = 0
$ws.Range("L792").Value = 0
$ws.Range("M792").Value = 4
$ws.Range("N792").Value = 0
$ws.Range("O792").Value = 0
$ws.Range("P792").Value = 0

# Row 793
$ws.Range("A793").NumberFormat = $dateFmt
$ws.Range("A793").Value = 45679
$ws.Range("B793").Value = 168.3500061035156
$ws.Range("C793").Value = 168.8999938964844
$ws.Range("D793").Value = 162.4799957275391
$ws.Range("E793").Value = 165.2700042724609
$ws.Range("F793").Value = 3319651
$ws.Range("G793").Value = 2025
$ws.Range("H793").Value = 1
$ws.Range("I793").Value = 22
$ws.Range("J793").Value = 0
$ws.Range("K793").Value = 0
$ws.Range("L793").Value = 0
$ws.Range("M793").Value = 4
$ws.Range("N793").Value = 0
$ws.Range("O793").Value = 0
$ws.Range("P793").Value = 0

# Row 794
$ws.Range("A794").NumberFormat = $dateFmt
$ws.Range("A794").Value = 45680
$ws.Range("B794").Value = 166
$ws.Range("C794").Value = 171
$ws.Range("D794").Value = 164.3899993896484
$ws.Range("E794").Value = 170.3099975585938
$ws.Range("F794").Value = 2858247
$ws.Range("G794").Value = 2025
$ws.Range("H794").Value = 1
$ws.Range("I794").Value = 23
$ws.Range("J794").Value = 0
$ws.Range("K794").Value = 0
$ws.Range("L794").Value = 0
$ws.Range("M794").Value = 4
$ws.Range("N794").Value = 0
$ws.Range("O794").Value = 0
$ws.Range("P794").Value = 0

# Row 795
$ws.Range("A795").NumberFormat = $dateFmt
$ws.Range("A795").Value = 45681
$ws.Range("B795").Value = 170.3000030517578
$ws.Range("C795").Value = 173.2599945068359
$ws.Range("D795").Value = 166.6100006103516
$ws.Range("E795").Value = 167.3699951171875
$ws.Range("F795").Value = 3192164
$ws.Range("G795").Value = 2025
$ws.Range("H795").Value = 1
$ws.Range("I795").Value = 24
$ws.Range("J795").Value = 0
$ws.Range("K795").Value = 0
$ws.Range("L795").Value = 0
$ws.Range("M795").Value = 4
$ws.Range("N795").Value = 0
$ws.Range("O795").Value = 0
$ws.Range("P795").Value = 0

# Row 796
$ws.Range("A796").NumberFormat = $dateFmt
$ws.Range("A796").Value = 45684
$ws.Range("B796").Value = 165.8999938964844
$ws.Range("C796").Value = 167.7299957275391
$ws.Range("D796").Value = 162.1999969482422
$ws.Range("E796").Value = 166.8800048828125
$ws.Range("F796").Value = 3403596
$ws.Range("G796").Value = 2025
$ws.Range("H796").Value = 1
$ws.Range("I796").Value = 27
$ws.Range("J796").Value = 0
$ws.Range("K796").Value = 0
$ws.Range("L796").Value = 0
$ws.Range("M796").Value = 5
$ws.Range("N796").Value = 0
$ws.Range("O796").Value = 0
$ws.Range("P796").Value = 0

# Row 797
$ws.Range("A797").NumberFormat = $dateFmt
$ws.Range("A797").Value = 45685
$ws.Range("B797").Value = 167.1999969482422
$ws.Range("C797").Value = 172.7299957275391
$ws.Range("D797").Value = 166.6600036621094
$ws.Range("E797").Value = 167.7799987792969
$ws.Range("F797").Value = 7808140
$ws.Range("G797").Value = 2025
$ws.Range("H797").Value = 1
$ws.Range("I797").Value = 28
$ws.Range("J797").Value = 0
$ws.Range("K797").Value = 0
$ws.Range("L797").Value = 0
$ws.Range("M797").Value = 5
$ws.Range("N797").Value = 0
$ws.Range("O797").Value = 0
$ws.Range("P797").Value = 0

# Row 798
$ws.Range("A798").NumberFormat = $dateFmt
$ws.Range("A798").Value = 45686
$ws.Range("B798").Value = 167.5099945068359
$ws.Range("C798").Value = 172.4900054931641
$ws.Range("D798").Value = 165.6100006103516
$ws.Range("E798").Value = 170.6000061035156
$ws.Range("F798").Value = 9475290
$ws.Range("G798").Value = 2025
$ws.Range("H798").Value = 1
$ws.Range("I798").Value = 29
$ws.Range("J798").Value = 0
$ws.Range("K798").Value = 0
$ws.Range("L798").Value = 0
$ws.Range("M798").Value = 5
$ws.Range("N798").Value = 0
$ws.Range("O798").Value = 0
$ws.Range("P798").Value = 0

# Row 799
$ws.Range("A799").NumberFormat = $dateFmt
$ws.Range("A799").Value = 45687
$ws.Range("B799").Value = 171.8300018310547
$ws.Range("C799").Value = 172
$ws.Range("D799").Value = 166.4700012207031
$ws.Range("E799").Value = 168.5899963378906
$ws.Range("F799").Value = 9873784
$ws.Range("G799").Value = 2025
$ws.Range("H799").Value = 1
$ws.Range("I799").Value = 30
$ws.Range("J799").Value = 0
$ws.Range("K799").Value = 0
$ws.Range("L799").Value = 0
$ws.Range("M799").Value = 5
$ws.Range("N799").Value = 0
$ws.Range("O799").Value = 0
$ws.Range("P799").Value = 0

